$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 73
$ws1.Range("F3").Value = 392
$ws1.Range("F4").Value = 444
$ws1.Range("F5").Value = 29
$ws1.Range("F6").Value = 14
$ws1.Range("F7").Value = 254
$ws1.Range("F8").Value = 14004
$ws1.Range("F9").Value = 107
$ws1.Range("F10").Value = 99
$ws1.Range("F11").Value = 5643
$ws1.Range("F12").Value = 579
$ws1.Range("F13").Value = 54
$ws1.Range("F15").Value = 52
$ws1.Range("F16").Value = 1226
$ws1.Range("F17").Value = 73
$ws1.Range("F18").Value = 164
$ws1.Range("F19").Value = 759
$ws1.Range("F20").Value = 2909
$ws1.Range("F21").Value = 46
$ws1.Range("F22").Value = 10428
$ws1.Range("F24").Value = 33
$ws1.Range("F25").Value = 50
$ws1.Range("F26").Value = 3708
$ws1.Range("F27").Value = 237

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 73
$ws4.Range("F3").Value = 392
$ws4.Range("F5").Value = 444
$ws4.Range("F6").Value = 29
$ws4.Range("F7").Value = 14
$ws4.Range("F8").Value = 254
$ws4.Range("F9").Value = 14004
$ws4.Range("F10").Value = 107
$ws4.Range("F11").Value = 99
$ws4.Range("F12").Value = 5643
$ws4.Range("F13").Value = 579
$ws4.Range("F14").Value = 54
$ws4.Range("F16").Value = 52
$ws4.Range("F17").Value = 1226
$ws4.Range("F18").Value = 73
$ws4.Range("F19").Value = 164
$ws4.Range("F20").Value = 759
$ws4.Range("F21").Value = 2909
$ws4.Range("F22").Value = 46
$ws4.Range("F24").Value = 10428
$ws4.Range("F26").Value = 33
$ws4.Range("F27").Value = 50
$ws4.Range("F28").Value = 3708
$ws4.Range("F29").Value = 237
